# Update the "取得日時" (retrieved timestamp) column on the "ランサーズ" sheet.
# Rows 2-9 had their timestamp bumped from 2025-11-23 01:38:33 to 2025-11-23 02:02:56
# as a result of a new scrape/append pass (commit: "Append: 2025-11-23 02:02 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-23 01:38:33"
$newTimestamp = "2025-11-23 02:02:56"

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
